$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("D4").Value = 45020
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 12000
$ws.Range("O4").Value = 12000
$ws.Range("P4").Value = 12000
$ws.Range("Q4").Value = '$/caja 18 kilos granel'
$ws.Range("S4").Value = 667

# Row 5
$ws.Range("D5").Value = 45021
$ws.Range("M5").Value = 50
$ws.Range("Q5").Value = '$/caja 18 kilos granel'

# Row 6
$ws.Range("D6").Value = 44699
$ws.Range("L6").Value = 'Especial'
$ws.Range("N6").Value = 13000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 13000
$ws.Range("Q6").Value = '$/caja 15 kilos granel'
$ws.Range("R6").Value = 'Provincia de Curicó'
$ws.Range("S6").Value = 867
$ws.Range("T6").Value = 15

# Row 7
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 120
$ws.Range("N7").Value = 11000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 11500
$ws.Range("S7").Value = 767

# Row 8
$ws.Range("D8").Value = 45040
$ws.Range("L8").Value = 'Especial'
$ws.Range("M8").Value = 50
$ws.Range("N8").Value = 13000
$ws.Range("O8").Value = 13000
$ws.Range("P8").Value = 13000
$ws.Range("Q8").Value = '$/caja 18 kilos empedrada'
$ws.Range("R8").Value = "Región de O'Higgins"
$ws.Range("S8").Value = 722
$ws.Range("T8").Value = 18

# Row 9
$ws.Range("D9").Value = 45040
$ws.Range("M9").Value = 40
$ws.Range("Q9").Value = '$/caja 18 kilos empedrada'
